# Append: 2025-12-02 06:29 JST
# Re-scrape update: the "ランサーズ" sheet rolls forward — rows 8-25 (the old
# tail of the listing) are dropped, and rows 2-7 are overwritten with the
# newest 6 listings (timestamp bumped to 06:29:37), mirroring a fresh
# scrape that only keeps 6 results this time.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1) Drop the old tail rows (8-25) entirely.
$ws.Range("A8:H25").EntireRow.Delete()

# 2) Clear out all existing hyperlink relationships so we can rebuild a
#    clean set for F2:F7 (mutating .Address in place would leave stale
#    duplicate relationships behind).
$ws.Hyperlinks.Delete()

# 3) Overwrite rows 2-7 with the new scrape results.

# Row 2
$ws.Range("A2").Value = "2025-12-02 06:29:37"
$ws.Range("B2").Value = "【せどり×ツール製作】APIを使用したせどりツールを製作できるエンジニアさんを募集します♪"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5217096"
$ws.Range("G2").Value = 243
$ws.Range("H2").Value = "🔥API ◆ツール"

# Row 3
$ws.Range("A3").Value = "2025-12-02 06:29:37"
$ws.Range("B3").Value = "管理システムの開発"
$ws.Range("C3").Value = "システム開発"
$ws.Range("D3").Value = "20,000 円 ~ 50,000 円 / 固定"
$ws.Range("E3").Value = "期限情報なし"
$ws.Range("F3").Value = "https://www.lancers.jp/work/detail/5445265"
$ws.Range("G3").Value = 103
$ws.Range("H3").Value = "◆開発 ◇管理"

# Row 4
$ws.Range("A4").Value = "2025-12-02 06:29:37"
$ws.Range("B4").Value = "【急募】iPhone用電子黒板アプリ開発者を探しています!"
$ws.Range("C4").Value = "システム開発"
$ws.Range("D4").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E4").Value = "期限情報なし"
$ws.Range("F4").Value = "https://www.lancers.jp/work/detail/5445417"
$ws.Range("G4").Value = 100
$ws.Range("H4").Value = "◆開発 ◇アプリ"

# Row 5
$ws.Range("A5").Value = "2025-12-02 06:29:37"
$ws.Range("B5").Value = "賃貸保証業の管理システム構築依頼"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "1,000,000 円 ~ 3,000,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5445528"
$ws.Range("G5").Value = 60
$ws.Range("H5").Value = "◇管理"

# Row 6
$ws.Range("A6").Value = "2025-12-02 06:29:37"
$ws.Range("B6").Value = "【UTAGE構築代行】各種初期設定・ステップ配信・会員サイトの作成など"
$ws.Range("C6").Value = "システム開発"
$ws.Range("D6").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E6").Value = "期限情報なし"
$ws.Range("F6").Value = "https://www.lancers.jp/work/detail/5445495"
$ws.Range("G6").Value = 38
$ws.Range("H6").Value = "◇サイト"

# Row 7 — this listing has no skill-tag summary, so H7 stays empty.
$ws.Range("A7").Value = "2025-12-02 06:29:37"
$ws.Range("B7").Value = "【若手歓迎×リモートOK】SRE/インフラエンジニア(Google Cloud/長期・金融系案件)"
$ws.Range("C7").Value = "システム開発"
$ws.Range("D7").Value = "500,000 円 ~ 1,000,000 円 / 固定"
$ws.Range("E7").Value = "期限情報なし"
$ws.Range("F7").Value = "https://www.lancers.jp/work/detail/5445466"
$ws.Range("G7").Value = 25
$ws.Range("H7").ClearContents()

# 4) Re-add the hyperlinks for the URL column, F2:F7, in row order so the
#    relationship ids come out sequential (rId1..rId6).
$ws.Hyperlinks.Add($ws.Range("F2"), "https://www.lancers.jp/work/detail/5217096", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5217096")
$ws.Hyperlinks.Add($ws.Range("F3"), "https://www.lancers.jp/work/detail/5445265", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5445265")
$ws.Hyperlinks.Add($ws.Range("F4"), "https://www.lancers.jp/work/detail/5445417", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5445417")
$ws.Hyperlinks.Add($ws.Range("F5"), "https://www.lancers.jp/work/detail/5445528", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5445528")
$ws.Hyperlinks.Add($ws.Range("F6"), "https://www.lancers.jp/work/detail/5445495", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5445495")
$ws.Hyperlinks.Add($ws.Range("F7"), "https://www.lancers.jp/work/detail/5445466", [Type]::Missing, [Type]::Missing, "https://www.lancers.jp/work/detail/5445466")

# 5) Column width tweaks (B, D, H). The COM ColumnWidth unit is offset by
#    -0.83 from the raw OOXML <col width> units on this host (verified by
#    round-tripping the sheet's existing columns), so subtract 0.83 from
#    the target raw widths (51, 32, 12) to land exactly on them.
$ws.Columns.Item(2).ColumnWidth = 50.17
$ws.Columns.Item(4).ColumnWidth = 31.17
$ws.Columns.Item(8).ColumnWidth = 11.17
